# Fixed spelling errors in ppt
#
# - Slide master + every slide layout: the cached "datetimeFigureOut"
#   field text changes from the 2-digit-year "11/4/14" to the
#   4-digit-year "11/4/2014".
# - Slide 2 title: "Eskrow" -> "Escrow".
# - Slide 3 body: "eskrow" -> "escrow" (both occurrences, each merged
#   with the neighbouring space like a normal find/replace would do).

$p = $ppt.ActivePresentation

function Fix-DatePlaceholder {
    param($shapes)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "11/4/14") {
                $tr.Text = "11/4/2014"
            }
        }
    }
}

# Slide master date placeholder.
Fix-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    Fix-DatePlaceholder $cl.Shapes
}

# --- Slide 2: "What is an Eskrow" -> "What is an Escrow" ---------------
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("Eskrow")
        if ($idx -ge 0) {
            $part = $tr.Characters($idx + 1, 6)
            $part.Text = "Escrow"
        }
    }
}

# --- Slide 3: fix both "eskrow" typos in the body placeholder ----------
$s3 = $p.Slides.Item(3)
$body = $s3.Shapes.Item(1)
$tr3 = $body.TextFrame.TextRange

# First occurrence: "Using the eskrow as a middle man ..."
# Replace "eskrow " (word + trailing space) with "escrow " so the space
# that used to start the next run is absorbed into the corrected word.
$first = $tr3.Characters(11, 7)
$first.Text = "escrow "

# Second occurrence: "Upon completion recipient gets money, the eskrow keeps 5%"
# Split "the " into its own run, then fix "eskrow " the same way.
$thePart = $tr3.Characters(128, 4)
$thePart.Text = "the "

$second = $tr3.Characters(132, 7)
$second.Text = "escrow "
